$d = $word.ActiveDocument

function Set-ParagraphXml($paragraph, $xmlPayload) {
    $beforeCount = $d.Paragraphs.Count
    $r = $paragraph.Range
    $r.End = $r.End - 1
    $r.Text = ""
    $r.InsertXML($xmlPayload)

    # Work around an edge case where replacing the content of the very
    # last paragraph in the body (the one immediately preceding
    # </w:body>'s sectPr) leaves a stray empty trailing paragraph behind.
    # If that happened, collapse the spurious paragraph mark back out.
    $afterCount = $d.Paragraphs.Count
    if ($afterCount -gt $beforeCount) {
        $extra = $afterCount - $beforeCount
        for ($i = 0; $i -lt $extra; $i++) {
            $lastP = $d.Paragraphs($d.Paragraphs.Count)
            if ($lastP.Range.Text.Length -le 1) {
                $fixRange = $d.Range($lastP.Range.Start - 1, $lastP.Range.End)
                $fixRange.Delete()
            }
        }
    }
}

$xmlPt = @'
<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:r><w:t xml:space="preserve">Sistemas de unidades: Definição das Unidades mais usadas em Engenharia e transformações entre sistemas. </w:t><w:br/><w:t xml:space="preserve">Estrutura atômica: Natureza elétrica da matéria. A carga do elétron. O núcleo do átomo. Teoria quântica: A radiação, os quanta e os fótons. Espectros de emissão e de absorção atômica. A dualidade onda-partícula da matéria. O princípio da incerteza. Os orbitais atômicos. Os números quânticos. Configuração eletrônica dos elementos. Partículas Elementares. </w:t><w:br/><w:t xml:space="preserve">Tabela periódica: A Lei e a tabela Periódica. Propriedades periódicas dos elementos, átomos e íons. </w:t><w:br/><w:t>Ligação Química: A ligação covalente. Estrutura de Lewis. Orbitais moleculares: Limitações da teoria de ligação de valência. Hibridização. Polaridade da ligação. Geometria molecular (Modelo VSEPR). Ligação Iônica. A classificação dos sólidos. As propriedades das ligações. Os compostos de coordenação. Complexos metálicos (teoria do campo cristalino). Ligação Metálica.</w:t><w:br/><w:br/><w:t>Nomenclatura de compostos inorgânicos: Funções Inorgânicas: ácidos; bases; sais; óxidos e nomenclaturas.</w:t><w:br/><w:t xml:space="preserve">Definições de ácidos e bases: Ácidos e bases (Arrhenius, Bronsted-Lowry e Lewis). </w:t><w:br/><w:t>Forças intermoleculares: Forças intermoleculares, líquidos e sólidos</w:t><w:br/><w:t>Soluções: Natureza das soluções. Dispersões coloidais e suspensões. Propriedades físicas e químicas. Tipos de soluções. Unidades e cálculos de concentração (Molaridade, fração molar, ppm, normalidade, molalidade, diluição). O processo de dissolução. Calor de dissolução. Solubilidade e temperatura.</w:t><w:br/><w:t xml:space="preserve">Gases (ideais e reais): Variáveis de estado. Lei combinada dos gases. Experiência de Torriceli. Pressão parcial dos gases. Teoria cinética dos gases. Gás ideal e real. Princípio de Avogadro. </w:t><w:br/><w:t>Reações químicas em solução aquosa: Principais reações químicas (ácido-base, precipitação, óxido-redução e complexação). Exemplos de reações formadoras de gases. Princípios de titulações ácido-base e de óxido-redução.</w:t><w:br/><w:t>Estequiometria e cálculos em química: Balanceamento de reações, cálculos estequiométricos, reagentes limitantes e rendimentos.</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>
'@

$xmlEn = @'
<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:r><w:rPr><w:i/></w:rPr><w:t>Unit systems: Definition of the Units most used in Engineering and transformations between systems.</w:t><w:br/><w:t>Atomic structure: Electric nature of matter. The charge of the electron. The nucleus of the atom. Quantum theory: Radiation, quanta and photons. Emission and atomic absorption spectra. The wave-particle duality of matter. The uncertainty principle. Atomic orbitals. Quantum numbers. Electronic configuration of the elements. Elementary Particles.</w:t><w:br/><w:t>Periodic table: The Law and the Periodic table. Periodic properties of elements, atoms and ions.</w:t><w:br/><w:t>Chemical Bond: The covalent bond. Lewis structure. Molecular orbitals: Limitations of the valence bond theory. Hybridization. Bond polarity. Molecular geometry (Model VSEPR). Ionic bonding. The classification of solids. The properties of the chemical bonds. Coordination compounds. Metal complexes (crystalline field theory). Metallic bond.</w:t><w:br/><w:t>Nomenclature of inorganic compounds: Inorganic Functions: acids; bases; salts; oxides and nomenclatures.</w:t><w:br/><w:t>Definitions of acids and bases: Acids and bases (Arrhenius, Bronsted-Lowry and Lewis).</w:t><w:br/><w:t>Intermolecular forces: Intermolecular forces, liquids and solids.</w:t><w:br/><w:t>Solutions: Nature of solutions. Colloidal dispersions and suspensions. Physical and chemical properties. Types of solutions. Concentration units and calculations (Molarity, molar fraction, ppm, normality, molality, dilution). The dissolution process. Heat of dissolution. Solubility and temperature.</w:t><w:br/><w:t>Gases (ideal and real): State variables. Combined gas law. Torriceli's Experience. Partial gas pressure. Kinetic theory of gases. Ideal and real gas. Avogadro's principle.</w:t><w:br/><w:t>Chemical reactions in aqueous solution: Main chemical reactions (acid-base, precipitation, oxide-reduction and complexation). Examples of gas-forming reactions. Principles of acid-base and oxide-reduction titrations.</w:t><w:br/><w:t>Stoichiometry and calculations in chemistry: Balancing reactions, stoichiometric calculations, limiting reagents and yields.</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>
'@

$xmlBib = @'
<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:r><w:t>ATKINS, Peter., Princípios de Química, questionando a vida moderna e o meio ambiente. 3ª Ed. Porto Alegre: Editora Bookman, 2006</w:t><w:br/><w:t>BRADY, J ; HUMISTON, G.E. Química geral. Rio de Janeiro: Ed. Livros Técnicos Científicos, 1981</w:t><w:br/><w:t>BROWN, T.L. ET al. Química a ciência central. 9.ed. São Paulo: Pearson Prentice Hall, 2005-2007</w:t><w:br/><w:t>CHANG, Raymond. Química geral: conceitos essenciais. 4.ed. s.l.:Ed. AMGH Editora Ltda., 2010.</w:t><w:br/><w:t>RUSSEL, J.B. Química geral. São Paulo: MacGrall-Hill</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>
'@

# Paragraph with the Portuguese "Programa" body text
$targetPt = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text.StartsWith("Sistemas de unidades: Defini")) {
        $targetPt = $p
        break
    }
}
Set-ParagraphXml $targetPt $xmlPt

# Paragraph with the English (italic) translation of the body text
$targetEn = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text.StartsWith("Unit systems: Definition of the Units")) {
        $targetEn = $p
        break
    }
}
Set-ParagraphXml $targetEn $xmlEn

# Paragraph with the Bibliografia references
$targetBib = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text.StartsWith("ATKINS, Peter")) {
        $targetBib = $p
        break
    }
}
Set-ParagraphXml $targetBib $xmlBib
